# Insert a new weekly record at row 609 for "Betarraga" (Hortaliza),
# shifting all subsequent rows down by one (old row 609 -> 610, ... old
# row 637 -> 638), matching the source diff which prepends a new
# observation dated 2023-08-09 (serial 45147) ahead of the existing series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 609; Excel copies the row-609
# formatting down onto the row that used to be 609 (now 610), and the new
# blank row 609 inherits formatting from the row above it (row 608).
$ws.Range("A609").EntireRow.Insert()

$ws.Range("A609").Value = 5
$ws.Range("B609").Value = "Macroferia Regional de Talca"
$ws.Range("C609").Value = "Maule"
$ws.Range("D609").Value = "8/9/2023"
$ws.Range("E609").Value = 7
$ws.Range("F609").Value = 100114014
$ws.Range("G609").Value = "Betarraga"
$ws.Range("H609").Value = "Sin especificar"
$ws.Range("I609").Value = "Primera"
$ws.Range("J609").Value = 5000
$ws.Range("K609").Value = 500
$ws.Range("L609").Value = 500
$ws.Range("M609").Value = 500
$ws.Range("N609").Value = "$/paquete 5 unidades"
$ws.Range("O609").Value = "Región del Maule"
$ws.Range("P609").Value = 100
$ws.Range("Q609").Value = 5
$ws.Range("R609").Value = "Hortaliza"
